$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header D1: "Context" -> "Theme"
$ws.Range("D1").Value = "Theme"

# New table data: row, A (paper id), B (which paper), C (title), D (theme)
$data = @(
  ,@(2, 8, "2018_2b", "2b. No money worries: new student identities for Graduate Level Apprentices", "graduate apprenticeships")
  ,@(3, 10, "2018_2c", "2c. The enemies within - inhibitors to learning", "policy setting")
  ,@(4, 57, "2019_3294016.3294025", "A Flexible Approach to Introductory Programming", "Teaching tech for engagement")
  ,@(5, 70, "2019_3294016.3298736", "The Institute of Coding: Addressing the UK Digital Skills Crisis", "a government-led initative")
  ,@(6, 82, "2020_3372356.3372362", "Computer Science Degree Accreditation in the UK:  A Post-Shadbolt Review Update", "accreditation")
  ,@(7, 86, "2020_3372356.3372364", "Errors and Poor Practices of Software Engineering  Students in Using Git", "professional tools")
  ,@(8, 87, "2020_3372356.3372364", "Errors and Poor Practices of Software Engineering  Students in Using Git", "misconceptions")
  ,@(9, 90, "2020_3372356.3372366", "Increasing academic diversity and inter-disciplinarity of  Computer Science in Higher Education", "widening participation (at a pinch could be `"recruitment/progression`")")
  ,@(10, 125, "2022_3498343.3498344", "Narrowing and Stretching: Addressing the Challenge of  Multi-track Programming", "Could be `"ability`". Could be `"pathways`". Neither fit well")
  ,@(11, 133, "2022_3498343.3498348", "Feedback and Engagement on an Introductory Programming  Module", "Could be `"assessment technique`" but it's a bit more subtle than that. Could be `"educational technology`", but that's misleading. Is engagement part of `"attitudes`"?")
  ,@(12, 135, "2022_3498343.3498349", "Co-Constructing a Community of Practice for Early-Career  Computer Science Academics in the UK", "professional  development")
  ,@(13, 168, "2024_3633053.3633057", "Incorporating Generative AI into Software Development Education", "generative AI")
  ,@(14, 185, "2025_3702212.3702214", "Learning without Limits: Analysing the Usage of Generative AI in  a Summative Assessment", "generative AI")
  ,@(15, 186, "2025_3702212.3702214", "Learning without Limits: Analysing the Usage of Generative AI in  a Summative Assessment", "generative AI")
  ,@(16, 187, "2025_3702212.3702215", "Group Assignments and Support Aimed to Develop Student Teamwork Skills and a Positive Attitude Towards Teamwork in  Computer Science Higher Education", "could be `"assessment techniques`". But there's nothing for peer learning/evaluation/feedback")
  ,@(17, 189, "2025_3702212.3702216", "Trunk and Branch: Fostering autonomous peer supportive learning environments through delivery & assessment", "_another one_ that's concerned with peer learning/feedback. Which is more than `"teaching/learning techniques`" I think.")
  ,@(18, 191, "2025_3702212.3702217", "Themes in the Declared Use of Generative Artificial Intelligence  in Assessment", "Could be `"assessment techniques`". More usefully Generative AI")
  ,@(19, 192, "2025_3702212.3702217", "Themes in the Declared Use of Generative Artificial Intelligence  in Assessment", "generative AI")
  ,@(20, 193, "2025_3702212.3702218", "Practical Insights for Engaging in Charity-University Collaborations for Computing Outreach for Disadvantaged  Young People", "working with charity. Disadvantage")
  ,@(21, 195, "2025_3702212.3702219", "FLARE: A Framework Supporting Code Comprehension and Formative Assessment in Block-Based Programming Education", "I could tick teaching/learning, but that would be misleading. It needs to be teacher development")
  ,@(22, 199, "2025_3702212.3702222", "Enhancing Learning and Teaching Experience for International  Students in Computing Subjects", "teaching/enhancement")
  ,@(23, 201, "2025_3702212.3702223", "Assessing Software Engineering Students’ Analytical Skills in the  Era of Generative AI", "generative AI? Analytic skills? ")
  ,@(24, 202, "2025_3702212.3702223", "Assessing Software Engineering Students’ Analytical Skills in the  Era of Generative AI", "generative AI")
  ,@(25, 203, "2025_3702212.3702224", "Where Have All the Papers Gone? Priming the pump of  pedagogical publishing in Europe", "publication? recognition of teaching? ")
)

foreach ($item in $data) {
  $row = $item[0]
  $a = $item[1]
  $b = $item[2]
  $c = $item[3]
  $d = $item[4]
  $ws.Cells.Item($row, 1).Value = $a
  $ws.Cells.Item($row, 2).Value = $b
  $ws.Cells.Item($row, 3).Value = $c
  $ws.Cells.Item($row, 4).Value = $d
}

# Apply the border/bold/center style (same as existing column A cells) to new rows 24 and 25 in column A
$ws.Range("A2").Copy()
$ws.Range("A24").PasteSpecial(-4122)
$ws.Range("A2").Copy()
$ws.Range("A25").PasteSpecial(-4122)

Write-Host "Edit complete"
